# Updates the "Out of PO" roster table with a refreshed player/position/team
# list (adds Ochai Agbaji / Toronto Raptors, drops Donte DiVincenzo, and
# reorders several rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$players = @(
    "Kelly Oubre Jr.",
    "Malik Beasley",
    "Ochai Agbaji",
    "Alperen Sengün",
    "Domantas Sabonis",
    "Victor Wembanyama",
    "Kristaps Porzingis",
    "Jaden McDaniels",
    "Dyson Daniels",
    "Michael Porter Jr.",
    "Carlton Carrington",
    "Bilal Coulibaly",
    "Josh Hart",
    "Kel'el Ware",
    "Donovan Mitchell",
    "Cam Thomas",
    "Andrew Wiggins",
    "De'Andre Hunter"
)

$positions = @(
    "SG,SF",
    "SG,SF",
    "SG,SF",
    "C",
    "C",
    "C",
    "PF,C",
    "SF,PF",
    "PG,SG,SF",
    "SF,PF",
    "PG,SG",
    "SG,SF",
    "SG,SF,PF",
    "PF,C",
    "PG,SG",
    "SG,SF",
    "SF,PF",
    "SF,PF"
)

$teams = @(
    "Philadelphia 76ers",
    "Detroit Pistons",
    "Toronto Raptors",
    "Houston Rockets",
    "Sacramento Kings",
    "San Antonio Spurs",
    "Boston Celtics",
    "Minnesota Timberwolves",
    "Atlanta Hawks",
    "Denver Nuggets",
    "Washington Wizards",
    "Washington Wizards",
    "New York Knicks",
    "Miami Heat",
    "Cleveland Cavaliers",
    "Brooklyn Nets",
    "Miami Heat",
    "Cleveland Cavaliers"
)

for ($i = 0; $i -lt $players.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $players[$i]
    $ws.Cells.Item($row, 2).Value = $positions[$i]
    $ws.Cells.Item($row, 3).Value = $teams[$i]
}
